# Update countries & provincias Spain
# - India (row 17): refreshed totals
# - Tailandia (row 67): refreshed totals
# - Seychelles / Montserrat (rows 205/206): the two countries swapped places
#   in the source feed, so the figures that used to sit on the "Seychelles"
#   row now belong to "Montserrat" and vice-versa (country name itself
#   stays tied to the row, only the numbers move).
# - Footer timestamp bumped from 06:34 to 07:04.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- India (row 17) ---
$ws.Range("B17").Value = 56409
$ws.Range("C17").Value = 58
$ws.Range("D17").Value = 16790
$ws.Range("E17").Value = 37729
$ws.Range("G17").Value = 1
$ws.Range("H17").Value = 1890

# --- Tailandia (row 67) ---
$ws.Range("B67").Value = 3000
$ws.Range("C67").Value = 8
$ws.Range("D67").Value = 2784
$ws.Range("E67").Value = 161

# --- Seychelles (row 205) / Montserrat (row 206) swap places ---
# Montserrat now reports first (row 205), Seychelles moves to row 206.
$ws.Range("A205").Value = "Montserrat"
$ws.Range("D205").Value = 7
$ws.Range("F205").Value = 1
$ws.Range("H205").Value = 1

$ws.Range("A206").Value = "Seychelles"
$ws.Range("D206").Value = 8
$ws.Range("F206").Value = 0
$ws.Range("H206").Value = 0

# --- Footer timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 8 de Mayo de 2020 a las 07:04"
